$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 243, pushing the existing rows 243:284 down to 244:285
$ws.Rows("243:243").Insert()

# Populate the newly inserted row 243 with the new record's data
$ws.Range("A243").Value = 10
$ws.Range("B243").Value = "Vega Modelo de Temuco"
$ws.Range("C243").Value = "La Araucanía"
$ws.Range("D243").Value = 44637
$ws.Range("D243").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E243").Value = 9
$ws.Range("F243").Value = 100114013
$ws.Range("G243").Value = "Zanahoria"
$ws.Range("H243").Value = "Sin especificar"
$ws.Range("I243").Value = "Primera"
$ws.Range("J243").Value = 100
$ws.Range("K243").Value = 6000
$ws.Range("L243").Value = 6000
$ws.Range("M243").Value = 6000
$ws.Range("N243").Value = "$/saco 25 kilos"
$ws.Range("O243").Value = "Región de La Araucanía"
$ws.Range("P243").Value = 240
$ws.Range("Q243").Value = 25
$ws.Range("R243").Value = "Hortaliza"
